# daily auto push: 2025-10-04 22:27 UTC
# Appends the next day's log entry (row 62) to Sheet1, right after the
# existing last row (61): date "2025/10/05", weekday "日", hour 4, rank 45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 62

# Column A holds dates stored as plain text (e.g. "2025/10/04" on the row
# above), not real Excel dates. A leading apostrophe forces text entry so
# "2025/10/05" isn't auto-converted to a date serial; resetting the style
# back to "Normal" afterwards drops the quote-prefix formatting Excel would
# otherwise remember, matching the unstyled cells used by every other data
# row in the sheet.
$ws.Range("A" + $newRow).Value = "'2025/10/05"
$ws.Range("A" + $newRow).Style = "Normal"

$ws.Range("B" + $newRow).Value = "日"
$ws.Range("C" + $newRow).Value = 4
$ws.Range("D" + $newRow).Value = 45
